$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.225.35"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.169.19"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.75"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.09"
$ws.Range("E6").Value = "  +5.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.166.39"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.25"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  +5.94%  "
$ws.Range("E13").Value = "  +17.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.23"
$ws.Range("E14").Value = "  +8.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.685.18"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.226.91"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("E17").Value = "  +7.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.169.00"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.40"
$ws.Range("E20").Value = "  +7.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.95"
$ws.Range("E21").Value = "  +6.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.16"
$ws.Range("E22").Value = "  +12.79%  "
$ws.Range("E23").Value = "  +8.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.88"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.87"
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +14.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.90"
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("E29").Value = "  +8.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.06"
$ws.Range("E30").Value = "  +6.73%  "
$ws.Range("E31").Value = "  +13.56%  "
$ws.Range("E32").Value = "  +7.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +11.97%  "
$ws.Range("E35").Value = "  +6.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.91"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "477.10"
$ws.Range("E37").Value = "  +6.81%  "
$ws.Range("E38").Value = "  +8.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +7.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0422"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.121.97"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("E43").Value = "  +6.02%  "
$ws.Range("E44").Value = "  +16.18%  "
$ws.Range("E45").Value = "  +11.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.28"
$ws.Range("E46").Value = "  +5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0586"
$ws.Range("E47").Value = "  +13.50%  "
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("E50").Value = "  +11.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.38"
$ws.Range("E51").Value = "  +6.12%  "
